$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.002727508544921875
$ws.Range("D2").Value = 0.001238822937011719

$ws.Range("C5").Value = -0.04925537109375
$ws.Range("D5").Value = -0.0650634765625

$ws.Range("D7").Value = -0.27880859375

$ws.Range("C8").Value = 0.001157760620117188
$ws.Range("D8").Value = -0.0177154541015625

$ws.Range("C9").Value = -0.06256103515625
$ws.Range("D9").Value = -0.145751953125
